# The <id> value for this table cell was previously split across three runs:
#   "<id>"  +  "p056v_1"  +  "</id>"
# Collapse them into a single run containing the full literal text, matching
# the newly-downloaded tc/tcn/tl content layout used elsewhere in the doc.
$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p056v_1</id>", $false, $true, $false, $false,
                         $false, $true, 1, $false, "<id>p056v_1</id>", 2)
